# Applies the changes described in the commit diff:
#  - sheet "links": add new column H "legend"; rows where carrier (D) was
#    "belastingen" become "indirecte_belastingen", and a matching "legend"
#    value of "indirecte_belastingen" is written into column H.
#  - sheet "carriers": remove the obsolete "belastingen"/"#3498db" row
#    (old row 2), shifting all following rows up by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "links" (sheet1)
# ---------------------------------------------------------------
$links = $wb.Worksheets.Item("links")

# New header for column H
$links.Range("H1").Value = "legend"

# Rows whose carrier value needs to change from "belastingen" to
# "indirecte_belastingen", and which also get a "legend" value.
$rowsToFix = @(2, 25, 26, 27, 28, 29, 30, 31, 32, 33, 34, 37)

foreach ($r in $rowsToFix) {
    $links.Cells.Item($r, 4).Value = "indirecte_belastingen"   # column D = carrier
    $links.Cells.Item($r, 8).Value = "indirecte_belastingen"   # column H = legend
}

# ---------------------------------------------------------------
# Sheet "carriers" (sheet3)
# ---------------------------------------------------------------
$carriers = $wb.Worksheets.Item("carriers")

# Remove the old row 2 ("belastingen" / "#3498db") entirely, shifting
# everything below it up by one row.
$carriers.Rows.Item(2).Delete()

Write-Output "done"
